$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, shifting existing rows 152:257 down to 153:258
$ws.Range("A152:R152").EntireRow.Insert()

# Populate the newly inserted row 152 with the new record's data
$fecha = Get-Date -Year 2022 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(152, 1).Value = 4
$ws.Cells.Item(152, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(152, 3).Value = "Los Lagos"
$ws.Cells.Item(152, 4).Value = $fecha
$ws.Cells.Item(152, 5).Value = 10
$ws.Cells.Item(152, 6).Value = 100112008
$ws.Cells.Item(152, 7).Value = "Coliflor"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 400
$ws.Cells.Item(152, 11).Value = 1500
$ws.Cells.Item(152, 12).Value = 1500
$ws.Cells.Item(152, 13).Value = 1500
$ws.Cells.Item(152, 14).Value = "`$/unidad"
$ws.Cells.Item(152, 15).Value = "Región Metropolitana"
$ws.Cells.Item(152, 16).Value = 1500
$ws.Cells.Item(152, 17).Value = 1
$ws.Cells.Item(152, 18).Value = "Hortaliza"
